$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-я Кольцевая 58")

# Columns B, C, D and F hold values that look like plain numbers (meter
# readings / phone-like IDs) but must stay stored as text, matching the
# source data. Force a text number format before writing so Excel keeps
# them as strings instead of auto-converting to numeric cells.
$ws.Range("B3:D4").NumberFormat = "@"
$ws.Range("F3:F4").NumberFormat = "@"

$ws.Cells.Item(3, 1).Value = "3-я Кольцевая 58, кв. 23"
$ws.Cells.Item(3, 2).Value = "33"
$ws.Cells.Item(3, 3).Value = "543"
$ws.Cells.Item(3, 4).Value = "334"
$ws.Cells.Item(3, 5).Value = "01.10.2023 в 17:11:29"
$ws.Cells.Item(3, 6).Value = "5734991862"

$ws.Cells.Item(4, 1).Value = "3-я Кольцевая 58, кв. 34"
$ws.Cells.Item(4, 2).Value = "23"
$ws.Cells.Item(4, 3).Value = "345"
$ws.Cells.Item(4, 4).Value = "244"
$ws.Cells.Item(4, 5).Value = "01.10.2023 в 20:46:23"
$ws.Cells.Item(4, 6).Value = "Alex Pol ID 128446192"
